$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 33 (pushing existing rows 33-120 down to 35-122)
$ws.Rows("33:34").Insert()

# New row 33
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44614
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100112031
$ws.Range("G33").Value = "Poroto verde"
$ws.Range("H33").Value = "Brío"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = 2000
$ws.Range("N33").Value = '$/kilo'
$ws.Range("O33").Value = "Región de La Araucanía"
$ws.Range("P33").Value = 2000
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"

# New row 34
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 44614
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 100112031
$ws.Range("G34").Value = "Poroto verde"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 50
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 2000
$ws.Range("N34").Value = '$/kilo'
$ws.Range("O34").Value = "Región de La Araucanía"
$ws.Range("P34").Value = 2000
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"
